$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.741.58"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.602.92"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.78"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.70"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").Value = "1.828.35"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "1.591.62"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.04"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "0.0₃0743"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "209.89"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.13"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.23"
$ws.Range("E22").Value = "  -5.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.06"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.67"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.09"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.38"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").Value = "1.289.41"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("E35").Value = "  +19.91%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.591"
$ws.Range("E37").Value = "  -4.47%  "
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.827"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.90"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").Value = "1.739.94"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.46"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0103"
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.17"
$ws.Range("E48").Value = "  +20.45%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.102"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0514"
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.57"
$ws.Range("E51").Value = "  +2.53%  "
